$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("B4").Value = 833000000.0
$ws.Range("C4").Value = 944000000.0
$ws.Range("D4").Value = 1007000000.0
$ws.Range("E4").Value = 1026000000.0
$ws.Range("F4").Value = 1088000000.0

# Row 13 - Accounts Payable
$ws.Range("B13").Value = 674000000.0
$ws.Range("C13").Value = 671000000.0
$ws.Range("D13").Value = 625000000.0
$ws.Range("E13").Value = 684000000.0
$ws.Range("F13").Value = 765000000.0

# Row 14 - Accrued Expenses
$ws.Range("C14").Value = 769800000.0

# Row 22 - Pension and Post-Retirement Liabilities
$ws.Range("C22").Value = 274800000.0

# Row 23 - Long Term Tax Liability (Deferred)
$ws.Range("B23").Value = 855000000.0
$ws.Range("C23").Value = 860000000.0
$ws.Range("D23").Value = 896000000.0
$ws.Range("E23").Value = 900000000.0
$ws.Range("F23").Value = 1076000000.0

# Row 38 - Shareholders Equity (Tangible)
$ws.Range("C38").Value = -1446300000.0

# Row 39 - Net Debt
$ws.Range("G39").Value = 1105900000.0

# Row 40 - Total Debt
$ws.Range("G40").Value = 1611400000.0
